# Append a new data row (row 3) to the "Random" sheet, mirroring the
# layout of the existing rows (Date, totalScore, ... , Method).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

$newRow = 3
$sourceRow = 2

# Copy the date cell's formatting (style index) from the row above so the
# new timestamp keeps the same date/time display format, without creating
# a duplicate style entry.
$ws.Cells.Item($sourceRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = 42605.888078703705   # Date
$ws.Cells.Item($newRow, 2).Value = 66                    # totalScore
$ws.Cells.Item($newRow, 3).Value = 0                     # posWordPercentage
$ws.Cells.Item($newRow, 4).Value = 0                     # negWordPercentage
$ws.Cells.Item($newRow, 5).Value = 0                     # posPhrasePercentage
$ws.Cells.Item($newRow, 6).Value = 0                     # negPhrasePercentage
$ws.Cells.Item($newRow, 7).Value = 0                     # ElapsedMs
$ws.Cells.Item($newRow, 8).Value = 0                     # wordCount
$ws.Cells.Item($newRow, 9).Value = 0                     # sentenceCount
$ws.Cells.Item($newRow, 10).Value = 0                    # posWordCount
$ws.Cells.Item($newRow, 11).Value = 0                    # negWordCount
$ws.Cells.Item($newRow, 12).Value = 0                    # positivePhraseCount
$ws.Cells.Item($newRow, 13).Value = 0                    # negativePhraseCount
$ws.Cells.Item($newRow, 14).Value = "Random"             # Method
